# Insert a new "pt_max" column before the existing "boson" column (column F),
# shifting the old F:K columns (and their data/widths) one place to the right
# to G:L, then populate the new column with its header and a constant value
# of 50 for every data row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift columns F:K -> G:L (carries formatting/column widths along, like
# Excel's own Insert Column command).
$ws.Columns("F").Insert() | Out-Null

# New header cell.
$ws.Range("F1").Value = "pt_max"

# New data values (pt_max = 50 for every row of data).
$ws.Range("F2:F10").Value = 50

# Match the author's final UI selection state.
$ws.Range("E16").Select() | Out-Null
